$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.408.92"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "2.326.59"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  -1.02%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  +7.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0800"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("D15").Value = "2.687.45"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").Value = "2.339.19"
$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.28%  "

$ws.Range("D18").Value = "43.285.06"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.42%  "

$ws.Range("E35").Value = "  -7.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0700"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.88%  "

$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("E40").Value = "  +0.61%  "

$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").Value = "1.993.30"
$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("E43").Value = "  +6.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.78%  "

$ws.Range("E46").Value = "  -3.01%  "

$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "2.554.18"
$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
